# Update countries & provincias Spain
# COVID-19 "paises" table refresh: new timestamp, updated case counts for
# several countries, a handful of countries re-ranked (swapped rows) after
# the refresh, a brand-new country row (Santo Tome y Principe) inserted in
# the existing range, and one additional row appended at the bottom
# (Sudan del Sur), growing the table from 214 to 215 data/header rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp (row 1)
$ws.Range('A1').Value = 'Datos actualizados a 6 de Abril de 2020 a las 16:22'

# Row 4 - Estados Unidos
$ws.Range('B4').Value = 337392
$ws.Range('C4').Value = 719
$ws.Range('E4').Value = 309763
$ws.Range('G4').Value = 36
$ws.Range('H4').Value = 9652

# Row 17 - Austria
$ws.Range('B17').Value = 12236
$ws.Range('C17').Value = 185
$ws.Range('E17').Value = 8553

# Row 19 - Brasil
$ws.Range('B19').Value = 11450
$ws.Range('C19').Value = 196
$ws.Range('E19').Value = 10832
$ws.Range('G19').Value = 5
$ws.Range('H19').Value = 491

# Row 25 - Noruega
$ws.Range('E25').Value = 5654
$ws.Range('G25').Value = 3
$ws.Range('H25').Value = 74

# Row 53 - Argentina
$ws.Range('D53').Value = 325
$ws.Range('E53').Value = 1181
$ws.Range('F53').Value = 94

# Row 82 - Libano
$ws.Range('D82').Value = 60
$ws.Range('E82').Value = 462

# Rows 102-104 re-ranked: Malta/Nigeria/Mauricio -> Mauricio/Malta/Nigeria
$ws.Range('A102').Value = 'Mauricio'
$ws.Range('B102').Value = 244
$ws.Range('C102').Value = 17
$ws.Range('D102').Value = 7
$ws.Range('E102').Value = 230
$ws.Range('F102').Value = 1
$ws.Range('H102').Value = 7

$ws.Range('A103').Value = 'Malta'
$ws.Range('B103').Value = 241
$ws.Range('C103').Value = 14
$ws.Range('D103').Value = 5
$ws.Range('E103').Value = 236
$ws.Range('F103').Value = 3
$ws.Range('H103').Value = 0

$ws.Range('A104').Value = 'Nigeria'
$ws.Range('B104').Value = 232
$ws.Range('D104').Value = 33
$ws.Range('E104').Value = 194
$ws.Range('F104').Value = 2
$ws.Range('H104').Value = 5

# Row 126 - Trinidad yTobago
$ws.Range('E126').Value = 96
$ws.Range('G126').Value = 1
$ws.Range('H126').Value = 8

# Row 157 - Tanzania
$ws.Range('B157').Value = 24
$ws.Range('C157').Value = 2
$ws.Range('E157').Value = 20

# Rows 164-165 swapped: Nueva Caledonia/Libia -> Libia/Nueva Caledonia
$ws.Range('A164').Value = 'Libia'
$ws.Range('D164').Value = 0
$ws.Range('H164').Value = 1

$ws.Range('A165').Value = 'Nueva Caledonia'
$ws.Range('D165').Value = 1
$ws.Range('H165').Value = 0

# Rows 171-172 swapped: Dominica/Fiyi -> Fiyi/Dominica
$ws.Range('A171').Value = 'Fiyi'
$ws.Range('C171').Value = 2

$ws.Range('A172').Value = 'Dominica'
$ws.Range('C172').Value = 0

# Rows 181-182 swapped: Seychelles/San Cristobal y Nieves -> San Cristobal y Nieves/Seychelles
$ws.Range('A181').Value = 'San Cristobal y Nieves'
$ws.Range('A182').Value = 'Seychelles'

# Rows 187-188 swapped: Nepal/Zimbabue -> Zimbabue/Nepal
$ws.Range('A187').Value = 'Zimbabue'
$ws.Range('D187').Value = 0
$ws.Range('H187').Value = 1

$ws.Range('A188').Value = 'Nepal'
$ws.Range('D188').Value = 1
$ws.Range('H188').Value = 0

# Rows 196-198 re-ranked: Botsuana/San Bartolome/Nicaragua -> San Bartolome/Nicaragua/Botsuana
$ws.Range('A196').Value = 'San Bartolome'
$ws.Range('D196').Value = 1
$ws.Range('H196').Value = 0

$ws.Range('A197').Value = 'Nicaragua'
$ws.Range('D197').Value = 0
$ws.Range('H197').Value = 1

$ws.Range('A198').Value = 'Botsuana'

# Rows 201-202 swapped: Belice/Islas Turcas y Caicos -> Islas Turcas y Caicos/Belice
$ws.Range('A201').Value = 'Islas Turcas y Caicos'
$ws.Range('F201').Value = 0
$ws.Range('G201').Value = 0

$ws.Range('A202').Value = 'Belice'
$ws.Range('F202').Value = 1
$ws.Range('G202').Value = 1

# Row 205 becomes the new country Santo Tome y Principe; Gambia/Burundi/
# Islas Virgenes Britanicas/Anguila all shift down one row (206-209)
$ws.Range('A205').Value = 'Santo Tome y Principe'
$ws.Range('C205').Value = 4
$ws.Range('D205').Value = 0
$ws.Range('E205').Value = 4
$ws.Range('H205').Value = 0

$ws.Range('A206').Value = 'Gambia'
$ws.Range('B206').Value = 4
$ws.Range('D206').Value = 2
$ws.Range('E206').Value = 1
$ws.Range('H206').Value = 1

$ws.Range('A207').Value = 'Burundi'

$ws.Range('A208').Value = 'Islas Virgenes Britanicas'

$ws.Range('A209').Value = 'Anguila'
$ws.Range('B209').Value = 3
$ws.Range('E209').Value = 3

# Row 212 becomes Bonaire, San Eustaquio y Saba (shifted from 209);
# rows 213-214 become San Pedro y Miquelon / Timor Oriental
$ws.Range('A212').Value = 'Bonaire, San Eustaquio y Saba'
$ws.Range('B212').Value = 2
$ws.Range('E212').Value = 2

$ws.Range('A213').Value = 'San Pedro y Miquelon'

$ws.Range('A214').Value = 'Timor Oriental'

# New row 215 appended: Sudan del Sur (was row 212, now its own new last row)
$ws.Range('A215').Value = 'Sudan del Sur'
$ws.Range('B215').Value = 1
$ws.Range('C215').Value = 0
$ws.Range('D215').Value = 0
$ws.Range('E215').Value = 1
$ws.Range('F215').Value = 0
$ws.Range('G215').Value = 0
$ws.Range('H215').Value = 0
